{"js": "// The document contains a single table, one column, one value per row\n// (a Renaissance/ZGC benchmark results table). This script updates the\n// specific cells whose values changed, matching the authoritative diff:\n//   row 0  (cell text \"100\")      -> \"0M\"\n//   row 1  (cell text \"0.03\")     -> \"0M\"\n//   row 2  (cell text \"686\")      -> \"0M\"\n//   row 3  (cell text \"108\")      -> \"410\"\n//   row 4  (cell text \"0.00003\")  -> \"0.00002\"\n//   row 6  (cell text \"0.00044\")  -> \"0.00030\"\n//   row 7  (cell text \"0.00080\")  -> \"0.00051\"\n//   row 11 (cell text \"0.01014\")  -> \"0.02587\"\n//   row 43 (multi-run tab-separated row, starts \"293\")  -> single run \"100\"\n//   row 44 (multi-run tab-separated row, starts \"8\")    -> single run \"0.03\"\n//   row 45 (multi-run tab-separated row, starts \"1\")    -> single run \"686\"\n//\n// Using range.insertText(text, \"Replace\") on the cell body's range keeps\n// the existing run formatting (rFonts/sz) and collapses any extra runs /\n// tab characters in a cell down to the single replacement run, which is\n// exactly what the target OOXML diff shows.\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\nconst updates = [\n  [0, \"0M\"],\n  [1, \"0M\"],\n  [2, \"0M\"],\n  [3, \"410\"],\n  [4, \"0.00002\"],\n  [6, \"0.00030\"],\n  [7, \"0.00051\"],\n  [11, \"0.02587\"],\n  [43, \"100\"],\n  [44, \"0.03\"],\n  [45, \"686\"],\n];\n\nfor (const [rowIndex, newText] of updates) {\n  const cell = table.getCell(rowIndex, 0);\n  const range = cell.body.getRange();\n  range.insertText(newText, \"Replace\");\n}\n\nawait context.sync();\n", "ps1": "# The document contains a single table, one column, one value per row\n# (a Renaissance/ZGC benchmark results table). This script updates the\n# specific cells whose values changed, matching the authoritative diff:\n#   row 1  (1-based) (cell text \"100\")      -> \"0M\"\n#   row 2  (1-based) (cell text \"0.03\")     -> \"0M\"\n#   row 3  (1-based) (cell text \"686\")      -> \"0M\"\n#   row 4  (1-based) (cell text \"108\")      -> \"410\"\n#   row 5  (1-based) (cell text \"0.00003\")  -> \"0.00002\"\n#   row 7  (1-based) (cell text \"0.00044\")  -> \"0.00030\"\n#   row 8  (1-based) (cell text \"0.00080\")  -> \"0.00051\"\n#   row 12 (1-based) (cell text \"0.01014\")  -> \"0.02587\"\n#   row 44 (1-based) (multi-run tab-separated row, starts \"293\") -> single run \"100\"\n#   row 45 (1-based) (multi-run tab-separated row, starts \"8\")   -> single run \"0.03\"\n#   row 46 (1-based) (multi-run tab-separated row, starts \"1\")   -> single run \"686\"\n#\n# Assigning Cell.Range.Text replaces the cell's contents (up to, but not\n# including, the end-of-cell marker) with a single run that reuses the\n# existing run formatting (rFonts/sz) of the cell, and collapses any extra\n# runs / tab characters down to the single replacement run -- exactly what\n# the target OOXML diff shows.\n\n$d = $word.ActiveDocument\n$tbl = $d.Tables.Item(1)\n\n$updates = @(\n    @{ Row = 1;  Text = \"0M\" },\n    @{ Row = 2;  Text = \"0M\" },\n    @{ Row = 3;  Text = \"0M\" },\n    @{ Row = 4;  Text = \"410\" },\n    @{ Row = 5;  Text = \"0.00002\" },\n    @{ Row = 7;  Text = \"0.00030\" },\n    @{ Row = 8;  Text = \"0.00051\" },\n    @{ Row = 12; Text = \"0.02587\" },\n    @{ Row = 44; Text = \"100\" },\n    @{ Row = 45; Text = \"0.03\" },\n    @{ Row = 46; Text = \"686\" }\n)\n\nforeach ($u in $updates) {\n    $cell = $tbl.Cell($u.Row, 1)\n    $cell.Range.Text = $u.Text\n}\n"}
